$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $rowA, $rowB, $firstCol, $lastCol) {
    $rangeA = $ws.Range($ws.Cells.Item($rowA, $firstCol), $ws.Cells.Item($rowA, $lastCol))
    $rangeB = $ws.Range($ws.Cells.Item($rowB, $firstCol), $ws.Cells.Item($rowB, $lastCol))

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

# Column B = 2, Column AD = 30
Swap-Rows $ws 129 130 2 30
Swap-Rows $ws 148 149 2 30
